$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add progress note for Tejomay Padole (row 7)
$ws.Range("B7").Value = "Flutter environment setup in android studio, started learning dart basics "

# Widen column B to fit new content (stored XML width rounds to nearest
# 1/6-character increment in this engine, so feed the char width that
# lands closest to the target stored width of 58.5546875)
$ws.Columns.Item(2).ColumnWidth = 57.6667

# Move selection to C2
$ws.Range("C2").Select()
